$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-06 14:30:48"
$ws.Range("E3").Value = "2026-02-06 14:30:50"
$ws.Range("E4").Value = "2026-02-06 14:30:53"
$ws.Range("E5").Value = "2026-02-06 14:30:55"
$ws.Range("E6").Value = "2026-02-06 14:30:58"
$ws.Range("E7").Value = "2026-02-06 14:31:00"
$ws.Range("E8").Value = "2026-02-06 14:31:03"
$ws.Range("E9").Value = "2026-02-06 14:31:06"
$ws.Range("E10").Value = "2026-02-06 14:31:08"
$ws.Range("E11").Value = "2026-02-06 14:31:10"
$ws.Range("E12").Value = "2026-02-06 14:31:13"
$ws.Range("E13").Value = "2026-02-06 14:31:15"
$ws.Range("E14").Value = "2026-02-06 14:31:17"
$ws.Range("E15").Value = "2026-02-06 14:31:20"
$ws.Range("E16").Value = "2026-02-06 14:31:22"
$ws.Range("E17").Value = "2026-02-06 14:31:25"
$ws.Range("E18").Value = "2026-02-06 14:31:27"
$ws.Range("E19").Value = "2026-02-06 14:31:30"
$ws.Range("E20").Value = "2026-02-06 14:31:32"
$ws.Range("E21").Value = "2026-02-06 14:31:35"
$ws.Range("E22").Value = "2026-02-06 14:31:37"
$ws.Range("E23").Value = "2026-02-06 14:31:39"
$ws.Range("E24").Value = "2026-02-06 14:31:42"
$ws.Range("E25").Value = "2026-02-06 14:31:44"
$ws.Range("E26").Value = "2026-02-06 14:31:47"
$ws.Range("E27").Value = "2026-02-06 14:31:49"
$ws.Range("E28").Value = "2026-02-06 14:31:52"
$ws.Range("E29").Value = "2026-02-06 14:31:54"
$ws.Range("E30").Value = "2026-02-06 14:31:57"
$ws.Range("E31").Value = "2026-02-06 14:31:59"
$ws.Range("E32").Value = "2026-02-06 14:32:01"
$ws.Range("E33").Value = "2026-02-06 14:32:04"
$ws.Range("E34").Value = "2026-02-06 14:32:06"
$ws.Range("E35").Value = "2026-02-06 14:32:09"
$ws.Range("E36").Value = "2026-02-06 14:32:11"
